$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.171.78"
$ws.Range("E2").Value = "  +0.79%  "

$ws.Range("D3").Value = "3.767.10"
$ws.Range("E3").Value = "  +1.16%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = "  -0.25%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.17"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +0.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.32"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +1.23%  "

$ws.Range("D7").Value = "3.766.22"
$ws.Range("E7").Value = "  +1.13%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.521"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  +0.60%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.160"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +0.82%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.41"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  -1.14%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.449"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -0.22%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000259"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  -0.76%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.12"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -0.07%  "

$ws.Range("D15").Value = "4.388.56"
$ws.Range("E15").Value = "  +0.90%  "

$ws.Range("D16").Value = "3.765.94"
$ws.Range("E16").Value = "  +0.86%  "

$ws.Range("D17").Value = "68.138.28"
$ws.Range("E17").Value = "  +0.73%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.89"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  -2.00%  "

$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.00"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  -0.20%  "

$ws.Range("B20").Value = "TRON"
$ws.Range("C20").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.112"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +0.66%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.70"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +0.25%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "465.62"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -0.29%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.697"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -0.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000151"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +13.54%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.96"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +1.53%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.18"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -0.21%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.86"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -1.21%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.08"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -0.55%  "

$ws.Range("E29").Value = "  +0.15%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.77"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +0.43%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.28"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -0.38%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "29.85"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +1.21%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.15"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -3.13%  "

$ws.Range("B34").Value = "Aptos"
$ws.Range("C34").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.14"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  +1.42%  "

$ws.Range("B35").Value = "Binance-PegBSC-USD"
$ws.Range("C35").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  +0.08%  "

$ws.Range("D36").Value = "3.719.98"
$ws.Range("E36").Value = "  +1.20%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.101"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -0.42%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.48"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +2.49%  "

$ws.Range("E39").Value = "  +0.22%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +1.24%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.77"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +0.66%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.998"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -0.32%  "

$ws.Range("E43").Value = "  -0.01%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "44.07"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +16.39%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.300"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -1.51%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "47.03"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +3.97%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.91"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +0.17%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.40"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  -1.35%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "145.12"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +1.48%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "388.06"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +0.36%  "

$ws.Range("D51").Value = "2.781.94"
$ws.Range("E51").Value = "  +3.98%  "
